$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KVO")

# Rename the "index" column header to "i"
$ws.Cells.Item(1, 1).Value = "i"

# The index column (column A) currently holds a 1-based row counter
# (1..501) for the data rows 2..503. Convert it to a 0-based counter
# (0..500) to match the rest of the zero-based series used elsewhere
# in the workbook.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $orig = $cell.Value2
    $cell.Value = $orig - 1
}

$ws.Columns.Item(1).ColumnWidth = 3.17
